$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) Append a new "Redshift" section at the very end of the document (after
#    the "A MASH network..." paragraph, before the sectPr).
# ---------------------------------------------------------------------------
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastP.Range.End, $lastP.Range.End)
$redshiftXml = $pkgOpen + `
  '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Redshift</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Redshift is another alternative to Arnold that renders scenes much faster.</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Instead of a Standin it uses a &#8220;proxy&#8221; which is basically the exact same thing.</w:t></w:r></w:p>' + `
  $pkgClose
$insertPoint.InsertXML($redshiftXml)

# ---------------------------------------------------------------------------
# 2) "A MASH network is a procedural network..." paragraph: merge the runs
#    that were split around the gramStart/gramEnd proofErr markers into one
#    run (keep the pPr/rPr bold marker that was already present).
# ---------------------------------------------------------------------------
$mashTextPara = $d.Paragraphs.Item(24)
$mashTextXml = $pkgOpen + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>A MASH network is a procedural network inside of Maya. It&#8217;s useful for instantiating objects and adding property nodes such as noise, movement, etc.</w:t></w:r></w:p>' + `
  $pkgClose
$mashTextPara.Range.InsertXML($mashTextXml)

# ---------------------------------------------------------------------------
# 3) "MASH" heading paragraph: drop the lastRenderedPageBreak marker (it
#    moves to the new NOTE bullet added under Standin).
# ---------------------------------------------------------------------------
$mashHeadingPara = $d.Paragraphs.Item(23)
$mashHeadingXml = $pkgOpen + `
  '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>MASH</w:t></w:r></w:p>' + `
  $pkgClose
$mashHeadingPara.Range.InsertXML($mashHeadingXml)

# ---------------------------------------------------------------------------
# 4) "Standin" section (paragraphs 19-22): drop the proofErr wrapping (text
#    unaffected), and add a new NOTE bullet (carrying the
#    lastRenderedPageBreak marker removed from the MASH heading above).
# ---------------------------------------------------------------------------
$standinStart = $d.Paragraphs.Item(19).Range.Start
$standinEnd = $d.Paragraphs.Item(22).Range.End
$standinRange = $d.Range($standinStart, $standinEnd)
$standinXml = $pkgOpen + `
  '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Standin</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>A standin is a file type (.ass) used by Arnold to help the scene run faster by hiding objects until a render is created.</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Standins are useful for things like grass when the scene could run very slow due to the number of polygons.</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>A standin can be used with a MASH network to instantiate the object thousands of times in a scene.</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>(NOTE: A duplicated standin might not show up properly in the render view. Make sure to use &#8220;duplicate special&#8221; (ctrl+shift+d) to duplicate a standin.)</w:t></w:r></w:p>' + `
  $pkgClose
$standinRange.InsertXML($standinXml)

# ---------------------------------------------------------------------------
# 5) "Sampling - ..." paragraph: merge the runs split around the gramStart/
#    gramEnd proofErr markers into one run.
# ---------------------------------------------------------------------------
$samplingPara = $d.Paragraphs.Item(11)
$samplingXml = $pkgOpen + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Sampling &#8211; Has settings regarding reducing noise in a rendered view. Each option reduces noise by increasing the sampling. Noise comes from many different things such as Motion blur, depth of field, diffuse, specular, etc. Increasing sampling for the wrong rays can make the render times increase without helping to remove the noise, so it&#8217;s important to adjust only change the settings that need it.</w:t></w:r></w:p>' + `
  $pkgClose
$samplingPara.Range.InsertXML($samplingXml)

# ---------------------------------------------------------------------------
# 6) "File Output - ..." paragraph: merge the runs split around the
#    spellStart/spellEnd/gramStart/gramEnd proofErr markers into one run.
# ---------------------------------------------------------------------------
$fileOutputPara = $d.Paragraphs.Item(9)
$fileOutputXml = $pkgOpen + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>File Output &#8211; Has settings regarding the images we render. If we want to use a frame range, make sure to change the Frame/Animation ext to something like &#8220;name_#.ext&#8221;</w:t></w:r></w:p>' + `
  $pkgClose
$fileOutputPara.Range.InsertXML($fileOutputXml)

# ---------------------------------------------------------------------------
# 7) "If using Arnold, ... RenderView." paragraph: merge the runs split
#    around the spellStart/spellEnd proofErr markers into one run.
# ---------------------------------------------------------------------------
$renderViewPara = $d.Paragraphs.Item(6)
$renderViewXml = $pkgOpen + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>If using Arnold, the Arnold window is accessed via Arnold &gt; Open Arnold RenderView.</w:t></w:r></w:p>' + `
  $pkgClose
$renderViewPara.Range.InsertXML($renderViewXml)

Write-Host "edits applied"
